# Updated cryptos list with GitHub Actions
# Applies the latest price/volume(1h) snapshot to the Sheet1 crypto table,
# including the re-ranking of Aave/BabyDogeCoin and Mantle/InjectiveProtocol.
# Numeric-looking price strings are written with a temporary Text number
# format so Excel keeps them as plain text (matching the source data),
# then the style is reset back to Normal to avoid leaving stray formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.533.65"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.444.82"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.109"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000174"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "2.890.31"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "62.533.44"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "2.446.89"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "320.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "647.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").Value = "2.567.51"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").Value = "0.0₃0957"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.364"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "149.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "153.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0292"
$ws.Range("E45").Value = "  -6.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.601"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0899"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "
